# 设备清单.xlsx update:
#  1. Insert a new row for "三通道直流电源" / "胡韬" above the current row 9
#     ("羽毛球机"), pushing everything from old row 9 down to row 10.
#  2. Update the "最后编辑人" (last editor) cell D2 from 董海鹏 to 胡韬.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 9 (current "羽毛球机" row) and fill it in.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "三通道直流电源"
$ws.Range("B9").Value = "胡韬"

# 2) Record that 胡韬 made the latest edit.
$ws.Range("D2").Value = "胡韬"

# Keep the active selection on the merged "last editor" column, matching
# where the editor was last working.
$ws.Range("D2:D22").Select()
